$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    right under the header, pushing all existing quarters down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# Copy the style from the row below (an untouched data row) onto the
# A-column cell of the newly inserted row so it keeps the bold/boxed
# look used by every other row in that column.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 13
$summary.Cells.Item(2, 4).Value = 0.73

# Renumber column A (0..5) for every data row now that a new one exists.
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Build the brand-new "2022-Q3" sheet by duplicating the existing
#    "2022-Q2" sheet (now at index 3, right after the new summary row
#    was added... no - sheet index is unaffected by row edits). This
#    keeps every sheet-level style (sheetPr/pageMargins/header &
#    column-A formatting) identical to its siblings.
# ---------------------------------------------------------------------
$templateQ2 = $wb.Worksheets.Item("2022-Q2")
$templateQ2.Copy($null, $summary)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template had 16 data rows (rows 2..17); 2022-Q3 only has 13
# (rows 2..14), so drop the trailing 3 rows.
$q3.Range("A15:H17").EntireRow.Delete()

# Data rows, columns A..H. Column A/H are genuine numbers; columns
# B..G must remain TEXT (e.g. leading zeros in fund codes, fixed
# decimal strings) exactly like every other quarter sheet.
$rows = @(
    @(0,  "003713", "英大睿盛灵活配置混合A",           "2.83", "93.65", "7.22", "0.2043", 3),
    @(1,  "003714", "英大睿盛灵活配置混合C",           "2.19", "93.65", "7.22", "0.1581", 3),
    @(2,  "162203", "泰达宏利稳定混合",                 "3.13", "91.16", "3.95", "0.1236", 8),
    @(3,  "001678", "英大国企改革主题股票",             "1.55", "93.30", "7.07", "0.1096", 3),
    @(4,  "001607", "英大策略优选混合A",                "0.57", "91.98", "7.84", "0.0447", 1),
    @(5,  "012522", "英大稳固增强核心一年持有混合C",     "1.24", "27.71", "2.12", "0.0263", 3),
    @(6,  "012521", "英大稳固增强核心一年持有混合A",     "0.75", "27.71", "2.12", "0.0159", 3),
    @(7,  "003447", "英大睿鑫灵活配置混合C",            "0.21", "92.71", "7.52", "0.0158", 7),
    @(8,  "007133", "嘉实长青竞争优势股票A",            "0.24", "90.21", "5.16", "0.0124", 6),
    @(9,  "002005", "工银新得利混合",                   "0.57", "27.58", "1.74", "0.0099", 4),
    @(10, "003446", "英大睿鑫灵活配置混合A",            "0.07", "92.71", "7.52", "0.0053", 7),
    @(11, "007134", "嘉实长青竞争优势股票C",            "0.04", "90.21", "5.16", "0.0021", 6),
    @(12, "001608", "英大策略优选混合C",                "0.02", "91.98", "7.84", "0.0016", 1)
)

# Stage the whole block far away from the live data, entering the
# numeric-looking text with a leading apostrophe so Excel treats it as
# text, then bulk copy/paste-values-only over the real cells: that
# moves the correct typed value across without carrying along the
# "number stored as text" quote-prefix style the apostrophe triggers
# on the staging cells.
$stageRow = 200
$r = $stageRow
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$stageRange = $q3.Range("A" + $stageRow + ":H" + ($stageRow + $rows.Count - 1))
$stageRange.Copy()
$q3.Range("A2:H14").PasteSpecial(-4163)
$stageRange.Clear()

$q3.Range("A1").Select()
